$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") - copy the existing header
# formatting (bold, centered, thin border) from H1 so the same style is reused.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new columns I (I0) and J (IF), rows 2-33
$data = @(
    @(11, 11),
    @(6, 7),
    @(6, 6),
    @(6, 6),
    @(6, 6),
    @(6, 7),
    @(7, 9),
    @(6, 7),
    @(4, 5),
    @(3, 4),
    @(9, 9),
    @(7, 7),
    @(6, 6),
    @(7, 8),
    @(5, 5),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(5, 6),
    @(11, 11),
    @(3, 4),
    @(8, 8),
    @(7, 7),
    @(6, 7),
    @(5, 7),
    @(4, 5),
    @(5, 5),
    @(8, 8),
    @(9, 9),
    @(5, 5),
    @(6, 6),
    @(9, 9)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
